$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9568985104560852
$ws.Range("B1").Value = 2.820319652557373
$ws.Range("C1").Value = 5.421573638916016
$ws.Range("D1").Value = 2.091215133666992
$ws.Range("E1").Value = 1.178155183792114
